# Updated work experience + education
# Remove the last row of the Education table (Bachelor of Arts and Education,
# English / Alexandria University / Grade: good), which corresponds to the
# "Sep 2002 - May 2006" entry.

$d = $word.ActiveDocument

# The "Education" table is the 4th table in the document
# (Experience, Tools & Technologies, Languages, Education).
$eduTable = $d.Tables(4)

# Delete the final row (the Bachelor's degree / Alexandria University entry).
$lastRow = $eduTable.Rows($eduTable.Rows.Count)
$lastRow.Delete()

# Bump the document-wide default line spacing slightly (274 -> 276 twentieths
# of a line, i.e. 13.7pt -> 13.8pt multiple spacing) as picked up by Word on
# resave of the "Normal" and "List Paragraph" styles.
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.LineSpacingRule = 5  # wdLineSpaceMultiple
$normalStyle.ParagraphFormat.LineSpacing = 13.8
